$p = $ppt.ActivePresentation
$m = $p.SlideMaster
Write-Host ($m.TextStyles | Get-Member | Out-String)
